# Applies the odds/snapshot-timestamp update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 1.25
$ws.Range("I2").Value = 1.42
$ws.Range("K2").Value = 14
$ws.Range("Q2").Value = 1.5
$ws.Range("BH2").Value = "2026-02-24 08:05:24"

# Row 3
$ws.Range("Q3").Value = 1.95
$ws.Range("BH3").Value = "2026-02-24 08:05:24"

# Row 4
$ws.Range("G4").Value = 3.8
$ws.Range("BH4").Value = "2026-02-24 08:05:24"

# Row 5
$ws.Range("G5").Value = 1.97
$ws.Range("H5").Value = 5.3
$ws.Range("I5").Value = 6.2
$ws.Range("P5").Value = 1.6
$ws.Range("Q5").Value = 2.14
$ws.Range("BH5").Value = "2026-02-24 08:05:24"

# Row 6
$ws.Range("BH6").Value = "2026-02-24 08:05:24"

# Row 7
$ws.Range("G7").Value = 2.34
$ws.Range("H7").Value = 3.15
$ws.Range("P7").Value = 1.78
$ws.Range("BH7").Value = "2026-02-24 08:05:24"

# Row 8
$ws.Range("H8").Value = 2.32
$ws.Range("I8").Value = 2.58
$ws.Range("K8").Value = 2.98
$ws.Range("P8").Value = 1.34
$ws.Range("Q8").Value = 3.45
$ws.Range("BH8").Value = "2026-02-24 08:05:24"

# Row 9
$ws.Range("F9").Value = 1.48
$ws.Range("G9").Value = 1.58
$ws.Range("H9").Value = 8.6
$ws.Range("K9").Value = 4.6
$ws.Range("P9").Value = 1.65
$ws.Range("Q9").Value = 2.04
$ws.Range("BH9").Value = "2026-02-24 08:05:24"

# Row 10
$ws.Range("BH10").Value = "2026-02-24 08:05:24"
